$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-17 (columns A-T), reflecting updated NATMI
# ligand-receptor pair statistics ("Natmi following Dr Hou advice").
$rows = @(
    @("ECs","Mmp2","Sdc2","ECs",3,1,2.408091666666667,7.224275,0.003423762576085583,0.003423762576085583,3,1,3.057109,9.171327,0.02694952608666365,0.02694952608666365,7.361798706991668,66.256188362925,0.00009226877885876117,0.00009226877885876115),
    @("ECs","Mmp2","Sdc2","FAPs",3,1,2.408091666666667,7.224275,0.003423762576085583,0.003423762576085583,3,1,89.02755999999999,267.08268,0.7848102735793893,0.7848102735793893,214.3865253396667,1929.478728057,0.002687004044008601,0.002687004044008601),
    @("ECs","Mmp2","Sdc2","M2",3,1,2.408091666666667,7.224275,0.003423762576085583,0.003423762576085583,2,0.6666666666666666,0.184005,0.5520149999999999,0.0016220709001794,0.0016220709001794,0.443100907125,3.987908164125,0.000005553585643791683,0.000005553585643791681),
    @("ECs","Mmp2","Sdc2","sCs",3,1,2.408091666666667,7.224275,0.003423762576085583,0.003423762576085583,3,1,21.16964733333333,63.508942,0.1866181294337677,0.1866181294337677,50.97845132967223,458.80606196705,0.0006389361675744295,0.0006389361675744293),
    @("FAPs","Mmp2","Sdc2","ECs",3,1,693.5675046666665,2080.702514,0.9860963763700007,0.9860963763700007,3,1,3.057109,9.171327,0.02694952608666365,0.02694952608666365,2120.311460624008,19082.80314561608,0.02657483001894783,0.02657483001894783),
    @("FAPs","Mmp2","Sdc2","FAPs",3,1,693.5675046666665,2080.702514,0.9860963763700007,0.9860963763700007,3,1,89.02755999999999,267.08268,0.7848102735793893,0.7848102735793893,61746.62263576193,555719.6037218574,0.7738985669145847,0.7738985669145847),
    @("FAPs","Mmp2","Sdc2","M2",3,1,693.5675046666665,2080.702514,0.9860963763700007,0.9860963763700007,2,0.6666666666666666,0.184005,0.5520149999999999,0.0016220709001794,0.0016220709001794,127.61988869619,1148.57899826571,0.001599518236882131,0.001599518236882131),
    @("FAPs","Mmp2","Sdc2","sCs",3,1,693.5675046666665,2080.702514,0.9860963763700007,0.9860963763700007,3,1,21.16964733333333,63.508942,0.1866181294337677,0.1866181294337677,14682.57947565335,132143.2152808802,0.1840234611995861,0.1840234611995861),
    @("M2","Mmp2","Sdc2","ECs",3,1,1.968565666666667,5.905697,0.002798855853950866,0.002798855853950866,3,1,3.057109,9.171327,0.02694952608666365,0.02694952608666365,6.018119816657666,54.163078349919,0.00007542783884886015,0.00007542783884886015),
    @("M2","Mmp2","Sdc2","FAPs",3,1,1.968565666666667,5.905697,0.002798855853950866,0.002798855853950866,3,1,89.02755999999999,267.08268,0.7848102735793893,0.7848102735793893,175.2565980031067,1577.30938202796,0.002196570828448454,0.002196570828448454),
    @("M2","Mmp2","Sdc2","M2",3,1,1.968565666666667,5.905697,0.002798855853950866,0.002798855853950866,2,0.6666666666666666,0.184005,0.5520149999999999,0.0016220709001794,0.0016220709001794,0.3622259254949999,3.260033329455,0.000004539942634490465,0.000004539942634490464),
    @("M2","Mmp2","Sdc2","sCs",3,1,1.968565666666667,5.905697,0.002798855853950866,0.002798855853950866,3,1,21.16964733333333,63.508942,0.1866181294337677,0.1866181294337677,41.67384091584156,375.064568242574,0.0005223172440190614,0.0005223172440190613),
    @("sCs","Mmp2","Sdc2","ECs",3,1,5.402408666666666,16.207226,0.007681005199962796,0.007681005199962796,3,1,3.057109,9.171327,0.02694952608666365,0.02694952608666365,16.51575215654466,148.641769408902,0.0002069994500081965,0.0002069994500081965),
    @("sCs","Mmp2","Sdc2","FAPs",3,1,5.402408666666666,16.207226,0.007681005199962796,0.007681005199962796,3,1,89.02755999999999,267.08268,0.7848102735793893,0.7848102735793893,480.9632617161866,4328.66935544568,0.006028131792347513,0.006028131792347513),
    @("sCs","Mmp2","Sdc2","M2",3,1,5.402408666666666,16.207226,0.007681005199962796,0.007681005199962796,2,0.6666666666666666,0.184005,0.5520149999999999,0.0016220709001794,0.0016220709001794,0.9940702067099998,8.946631860389997,0.0000124591350189863,0.0000124591350189863),
    @("sCs","Mmp2","Sdc2","sCs",3,1,5.402408666666666,16.207226,0.007681005199962796,0.007681005199962796,3,1,21.16964733333333,63.508942,0.1866181294337677,0.1866181294337677,114.3670862238769,1029.303776014892,0.0014334148225881,0.0014334148225881)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowVals = $rows[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}

Write-Host "Updated $($rows.Count) data rows."
